$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Generic cleanup: wipe both contents and formatting for the stale
# "computed/derived" columns instead of hand-maintaining each value.
$ranges = @("C2:D4", "G2:L4")

foreach ($addr in $ranges) {
    $rng = $ws.Range($addr)
    $rng.ClearContents()
    $rng.ClearFormats()
}
